$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.461.27'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').Value = '1.726.18'
$ws.Range('E3').Value = '  -0.50%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9969'
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.91'
$ws.Range('E5').Value = '  -1.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9974'
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4893'
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2603'
$ws.Range('E8').Value = '  -2.45%  '
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('D10').Value = '1.719.51'
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06991'
$ws.Range('E11').Value = '  -0.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.60'
$ws.Range('E12').Value = '  -0.68%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.523'
$ws.Range('E14').Value = '  -1.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.20'
$ws.Range('E15').Value = '  -0.34%  '
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D17').Value = '26.456.61'
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9969'
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007169'
$ws.Range('E19').Value = '  -1.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.33'
$ws.Range('E20').Value = '  -1.75%  '
$ws.Range('D21').Value = '1.946.71'
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.462'
$ws.Range('E22').Value = '  -1.97%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.518'
$ws.Range('E23').Value = '  -2.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.113'
$ws.Range('E24').Value = '  -2.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '137.18'
$ws.Range('E25').Value = '  -2.30%  '
$ws.Range('E26').Value = '  -1.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.411'
$ws.Range('E27').Value = '  -0.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '106.70'
$ws.Range('E28').Value = '  -1.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.744'
$ws.Range('E29').Value = '  -2.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.907'
$ws.Range('E30').Value = '  -2.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08020'
$ws.Range('E31').Value = '  -0.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.646'
$ws.Range('E32').Value = '  -1.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04495'
$ws.Range('E33').Value = '  -1.61%  '
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.599'
$ws.Range('E35').Value = '  -0.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.001'
$ws.Range('E36').Value = '  -1.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6251'
$ws.Range('E37').Value = '  -2.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9396'
$ws.Range('E38').Value = '  +3.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.988'
$ws.Range('E39').Value = '  -2.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.384'
$ws.Range('E40').Value = '  -0.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9969'
$ws.Range('E41').Value = '  -0.70%  '
$ws.Range('E42').Value = '  -1.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.79'
$ws.Range('E43').Value = '  -1.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.383'
$ws.Range('E44').Value = '  -0.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3845'
$ws.Range('E45').Value = '  -1.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.895'
$ws.Range('E46').Value = '  -0.90%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1160'
$ws.Range('E47').Value = '  -2.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05367'
$ws.Range('E48').Value = '  -0.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.42'
$ws.Range('E49').Value = '  -0.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.730'
$ws.Range('E50').Value = '  -0.98%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.34'
$ws.Range('E51').Value = '  -0.49%  '
